$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''51.956.76'
$ws.Range("E2").Value = '  -1.00%  '

$ws.Range("D3").Value = '''2.928.31'
$ws.Range("E3").Value = '  +0.00%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").Value = '''358.03'
$ws.Range("E5").Value = '  +1.37%  '

$ws.Range("D6").Value = '''110.54'
$ws.Range("E6").Value = '  -2.19%  '

$ws.Range("E7").Value = '  +1.17%  '

$ws.Range("D9").Value = '''0.633'
$ws.Range("E9").Value = '  +0.59%  '

$ws.Range("D10").Value = '''39.43'
$ws.Range("E10").Value = '  -1.96%  '

$ws.Range("D11").Value = '''0.0879'
$ws.Range("E11").Value = '  +1.57%  '

$ws.Range("E12").Value = '  +0.95%  '

$ws.Range("D13").Value = '''19.69'
$ws.Range("E13").Value = '  -2.51%  '

$ws.Range("D14").Value = '''7.96'
$ws.Range("E14").Value = '  +1.28%  '

$ws.Range("D15").Value = '''3.392.20'
$ws.Range("E15").Value = '  +0.16%  '

$ws.Range("D16").Value = '''2.929.37'
$ws.Range("E16").Value = '  -0.42%  '

$ws.Range("E17").Value = '  -0.76%  '

$ws.Range("D18").Value = '''51.970.96'
$ws.Range("E18").Value = '  -0.99%  '

$ws.Range("E19").Value = '  +1.09%  '

$ws.Range("D20").Value = '''7.61'
$ws.Range("E20").Value = '  -1.67%  '

$ws.Range("D21").Value = '''14.09'
$ws.Range("E21").Value = '  -2.45%  '

$ws.Range("D22").Value = '''0.0₃0984'
$ws.Range("E22").Value = '  -0.05%  '

$ws.Range("E23").Value = '  -0.12%  '

$ws.Range("D24").Value = '''270.43'
$ws.Range("E24").Value = '  -0.46%  '

$ws.Range("E25").Value = '  +0.83%  '

$ws.Range("D26").Value = '''0.186'
$ws.Range("E26").Value = '  +12.43%  '

$ws.Range("D27").Value = '''27.14'
$ws.Range("E27").Value = '  +0.28%  '

$ws.Range("D28").Value = '''7.53'
$ws.Range("E28").Value = '  +15.23%  '

$ws.Range("E29").Value = '  +0.03%  '

$ws.Range("E30").Value = '  +13.97%  '

$ws.Range("D31").Value = '''10.62'
$ws.Range("E31").Value = '  -0.31%  '

$ws.Range("D32").Value = '''38.23'
$ws.Range("E32").Value = '  -0.01%  '

$ws.Range("D33").Value = '''2.29'
$ws.Range("E33").Value = '  +1.84%  '

$ws.Range("D34").Value = '''6.10'
$ws.Range("E34").Value = '  -1.91%  '

$ws.Range("D35").Value = '''52.27'
$ws.Range("E35").Value = '  -1.92%  '

$ws.Range("D36").Value = '''0.0445'
$ws.Range("E36").Value = '  -2.00%  '

$ws.Range("E37").Value = '  +0.06%  '

$ws.Range("D38").Value = '''3.26'
$ws.Range("E38").Value = '  -2.89%  '

$ws.Range("D39").Value = '''18.48'
$ws.Range("E39").Value = '  -1.70%  '

$ws.Range("D40").Value = '''2.01'
$ws.Range("E40").Value = '  -3.61%  '

$ws.Range("E41").Value = '  +0.33%  '

$ws.Range("E42").Value = '  +2.31%  '

$ws.Range("D43").Value = '''23.26'
$ws.Range("E43").Value = '  -5.59%  '

$ws.Range("D44").Value = '''119.42'
$ws.Range("E44").Value = '  -2.89%  '

$ws.Range("E45").Value = '  -1.63%  '

$ws.Range("E46").Value = '  -2.28%  '

$ws.Range("E47").Value = '  -4.90%  '

$ws.Range("D48").Value = '''2.137.17'
$ws.Range("E48").Value = '  -3.78%  '

$ws.Range("D49").Value = '''0.250'
$ws.Range("E49").Value = '  -5.12%  '

$ws.Range("D50").Value = '''0.0336'
$ws.Range("E50").Value = '  -1.04%  '

$ws.Range("D51").Value = '''9.17'
$ws.Range("E51").Value = '  -0.18%  '
